# Update odds values in the "Jogos do Dia" sheet to reflect the latest
# Betfair Back/Lay quotes, as captured by the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 3.15
$ws.Range("Q2").Value = 2.44

# Row 3
$ws.Range("P3").Value = 1.67
$ws.Range("Q3").Value = 1.91

# Row 4
$ws.Range("F4").Value = 2.32
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3.45
$ws.Range("J4").Value = 2.66
$ws.Range("K4").Value = 3.6
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.52

# Row 7
$ws.Range("F7").Value = 2.54
$ws.Range("G7").Value = 2.62
$ws.Range("H7").Value = 2.68
$ws.Range("I7").Value = 2.92
$ws.Range("J7").Value = 3.6
$ws.Range("K7").Value = 4.1
$ws.Range("P7").Value = 2.08
$ws.Range("Q7").Value = 1.78

# Row 8
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 4.6
$ws.Range("H8").Value = 1.99
$ws.Range("I8").Value = 2.12
$ws.Range("K8").Value = 3.75
$ws.Range("Q8").Value = 2.04

# Row 9
$ws.Range("H9").Value = 1.93

# Row 11
$ws.Range("F11").Value = 2.58
$ws.Range("G11").Value = 2.96
$ws.Range("H11").Value = 2.76
$ws.Range("J11").Value = 3.35
$ws.Range("K11").Value = 4
$ws.Range("P11").Value = 1.77

# Row 12
$ws.Range("H12").Value = 6
$ws.Range("I12").Value = 8.199999999999999
$ws.Range("J12").Value = 3.65
$ws.Range("P12").Value = 1.7
$ws.Range("Q12").Value = 2.14

# Row 13
$ws.Range("F13").Value = 1.86
$ws.Range("G13").Value = 2.28
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 4.8
$ws.Range("J13").Value = 3.3
$ws.Range("K13").Value = 5.1
$ws.Range("P13").Value = 1.96
$ws.Range("Q13").Value = 1.61
$ws.Range("T13").Value = 1.66
$ws.Range("U13").Value = 2.26
$ws.Range("X13").Value = 23
$ws.Range("AA13").Value = 95
$ws.Range("AO13").Value = 48

# Row 14
$ws.Range("F14").Value = 1.85
$ws.Range("G14").Value = 1.99
$ws.Range("H14").Value = 4.2
$ws.Range("I14").Value = 4.8
$ws.Range("J14").Value = 3.7
$ws.Range("K14").Value = 4.2
$ws.Range("Q14").Value = 1.83

# Row 16
$ws.Range("F16").Value = 4.8
$ws.Range("G16").Value = 5.6
$ws.Range("H16").Value = 1.83
$ws.Range("I16").Value = 1.94
$ws.Range("J16").Value = 3.4
$ws.Range("K16").Value = 3.85
$ws.Range("P16").Value = 1.8
$ws.Range("Q16").Value = 1.9
